$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.058215421723377
$ws.Range("D2").Value = 1.060164991995332
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.068920104755342
$ws.Range("I2").Value = 1.048232875903794
$ws.Range("J2").Value = 1.063207800014802
$ws.Range("K2").Value = 1.062892515171501
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.071623997119951
$ws.Range("N2").Value = 1.064717677042463

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.059368513081016
$ws.Range("D3").Value = 1.061071286016051
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.070023045119048
$ws.Range("I3").Value = 1.048564704089415
$ws.Range("J3").Value = 1.064012690452243
$ws.Range("K3").Value = 1.06361294887268
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.072542296131228
$ws.Range("N3").Value = 1.065523710516647

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.060114374261587
$ws.Range("D4").Value = 1.061657498226818
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.070736901240841
$ws.Range("I4").Value = 1.048778113773388
$ws.Range("J4").Value = 1.064532674355968
$ws.Range("K4").Value = 1.064078263646597
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.073136081476456
$ws.Range("N4").Value = 1.066044432857158

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.06042787187381
$ws.Range("D5").Value = 1.06190388933174
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.071037050018372
$ws.Range("I5").Value = 1.048867518947589
$ws.Range("J5").Value = 1.064751076538047
$ws.Range("K5").Value = 1.064273677965224
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.07338560990007
$ws.Range("N5").Value = 1.066263145195386

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.060480505878222
$ws.Range("D6").Value = 1.061945256412261
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.071087448896238
$ws.Range("I6").Value = 1.048882512165953
$ws.Range("J6").Value = 1.064787735570931
$ws.Range("K6").Value = 1.064306476938528
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.073427501034459
$ws.Range("N6").Value = 1.066299856288302

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.060118563475741
$ws.Range("D7").Value = 1.061660790722759
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.070740911671197
$ws.Range("I7").Value = 1.048779309636352
$ws.Range("J7").Value = 1.064535593437978
$ws.Range("K7").Value = 1.064080875583371
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.073139416075539
$ws.Range("N7").Value = 1.066047356084599

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.058605170189261
$ws.Range("D8").Value = 1.060471324171432
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.069292811705593
$ws.Range("I8").Value = 1.04834528921941
$ws.Range("J8").Value = 1.063479989570563
$ws.Range("K8").Value = 1.063136166205635
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.07193442693265
$ws.Range("N8").Value = 1.064990253138613

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.055936283564102
$ws.Range("D9").Value = 1.05837362797139
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.066742426153607
$ws.Range("I9").Value = 1.0475704812543
$ws.Range("J9").Value = 1.061613464903678
$ws.Range("K9").Value = 1.061464906840641
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.069807874992315
$ws.Range("N9").Value = 1.063121077792608

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.054155535678641
$ws.Range("D10").Value = 1.056973994920275
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.065043037580692
$ws.Range("I10").Value = 1.047047196409996
$ws.Range("J10").Value = 1.060364759768415
$ws.Range("K10").Value = 1.060346292399623
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.068387980110353
$ws.Range("N10").Value = 1.061870599352824

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.053384079654307
$ws.Range("D11").Value = 1.056367655041703
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.064307379639102
$ws.Range("I11").Value = 1.046819003947641
$ws.Range("J11").Value = 1.059823013558304
$ws.Range("K11").Value = 1.059860858992786
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.067772620498326
$ws.Range("N11").Value = 1.061328083800956

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.053097467629853
$ws.Range("D12").Value = 1.0561423892858
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.064034150884307
$ws.Range("I12").Value = 1.046734001290341
$ws.Range("J12").Value = 1.059621626379177
$ws.Range("K12").Value = 1.059680386515573
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.067543967153392
$ws.Range("N12").Value = 1.061126410628935

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.053158949575887
$ws.Range("D13").Value = 1.056190711562142
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.064092758157132
$ws.Range("I13").Value = 1.046752245598534
$ws.Range("J13").Value = 1.059664831801461
$ws.Range("K13").Value = 1.059719105767452
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.067593017769356
$ws.Range("N13").Value = 1.061169677407873

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.053360389426797
$ws.Range("D14").Value = 1.056349035407891
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.064284793912449
$ws.Range("I14").Value = 1.046811982536956
$ws.Range("J14").Value = 1.05980637007053
$ws.Range("K14").Value = 1.059845944368072
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.067753721606941
$ws.Range("N14").Value = 1.061311416677519

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.053484495308999
$ws.Range("D15").Value = 1.056446578088373
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.064403117075843
$ws.Range("I15").Value = 1.046848756380317
$ws.Range("J15").Value = 1.059893555435617
$ws.Range("K15").Value = 1.059924072473487
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.067852725733914
$ws.Range("N15").Value = 1.061398725855826

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.054206726509412
$ws.Range("D16").Value = 1.057014229547405
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.065091864665255
$ws.Range("I16").Value = 1.047062306904577
$ws.Range("J16").Value = 1.06040069150383
$ws.Range("K16").Value = 1.060378486493699
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.068428808137071
$ws.Range("N16").Value = 1.061906582115425

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.054659659971102
$ws.Range("D17").Value = 1.057370224352424
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.065523947626159
$ws.Range("I17").Value = 1.047195831050075
$ws.Range("J17").Value = 1.06071852312641
$ws.Range("K17").Value = 1.060663242486846
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.068790025340848
$ws.Range("N17").Value = 1.062224865095365

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.054923811807103
$ws.Range("D18").Value = 1.057577842206114
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.065775992561778
$ws.Range("I18").Value = 1.047273558433609
$ws.Range("J18").Value = 1.060903807820523
$ws.Range("K18").Value = 1.060829233039293
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.069000665591795
$ws.Range("N18").Value = 1.062410412914996

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.055013874603949
$ws.Range("D19").Value = 1.057648629816798
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.065861936496981
$ws.Range("I19").Value = 1.047300035186634
$ws.Range("J19").Value = 1.060966967991737
$ws.Range("K19").Value = 1.060885814113959
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.06907247974329
$ws.Range("N19").Value = 1.062473662780898

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.05461106829417
$ws.Range("D20").Value = 1.057332032368439
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.065477587332231
$ws.Range("I20").Value = 1.047181521207453
$ws.Range("J20").Value = 1.060684433266751
$ws.Range("K20").Value = 1.06063270150207
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.0687512755003
$ws.Range("N20").Value = 1.062190726824195

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.053301072051685
$ws.Range("D21").Value = 1.056302414204649
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.064228243421365
$ws.Range("I21").Value = 1.046794398176288
$ws.Range("J21").Value = 1.059764694961628
$ws.Range("K21").Value = 1.059808598009078
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.067706400586065
$ws.Range("N21").Value = 1.061269682385183

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.052477083728441
$ws.Range("D22").Value = 1.055654796775969
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.063442889134456
$ws.Range("I22").Value = 1.046549599271685
$ws.Range("J22").Value = 1.059185501382196
$ws.Range("K22").Value = 1.059289520168777
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.067048974936382
$ws.Range("N22").Value = 1.060689666284435

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.052913928503069
$ws.Range("D23").Value = 1.055998135469605
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.063859205665081
$ws.Range("I23").Value = 1.046679504583804
$ws.Range("J23").Value = 1.059492630159616
$ws.Range("K23").Value = 1.059564781676247
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.067397533687508
$ws.Range("N23").Value = 1.060997231219946

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.054633024893936
$ws.Range("D24").Value = 1.057349289767977
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.065498535496767
$ws.Range("I24").Value = 1.047187987688115
$ws.Range("J24").Value = 1.060699837319059
$ws.Range("K24").Value = 1.060646501975841
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.068768785044503
$ws.Range("N24").Value = 1.062206152752024

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.056626511354596
$ws.Range("D25").Value = 1.058916137038395
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.067401605913475
$ws.Range("I25").Value = 1.047771975498968
$ws.Range("J25").Value = 1.062096770861628
$ws.Range("K25").Value = 1.061897748053181
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.07035802353087
$ws.Range("N25").Value = 1.063605070100454

